{"js": "// 1) \"...but there's a problem\" + \":\" were two separate runs with the same\n//    formatting; merge them into a single run reading \"...problem:\" (no\n//    visible text change, only a run-structure normalization).\nconst body = context.document.body;\n\nconst colonSearch = body.search(\"but there\\u2019s a problem:\", { matchCase: true });\ncolonSearch.load(\"text\");\nawait context.sync();\n\nconst colonRange = colonSearch.items[0];\ncolonRange.insertText(colonRange.text, \"Replace\");\nawait context.sync();\n\n// 2) Insert a new clause right after \"...You bring no silver at all\" and\n//    before the following \". How would you reach the temple...\" sentence.\nconst anchorSearch = body.search(\"You bring no silver at all\", { matchCase: true });\nanchorSearch.load(\"text\");\nawait context.sync();\n\nconst anchorRange = anchorSearch.items[0];\nconst insertedText = \", and you can\\u2019t walk the path you\\u2019ve trodden, though you can cross path\";\nanchorRange.insertText(insertedText, \"After\");\nawait context.sync();\n\n// Force the newly inserted clause into its own run (distinct from the\n// surrounding text it was merged with) by toggling a character property\n// on it and back off again.\nconst newSearch = body.search(insertedText, { matchCase: true });\nnewSearch.load(\"text\");\nawait context.sync();\n\nconst newRange = newSearch.items[0];\nnewRange.font.bold = true;\nawait context.sync();\nnewRange.font.bold = false;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"...but there's a problem\" and the following \":\" are two separate runs\n#    sharing identical formatting; replacing the combined text with itself\n#    (via Find/Replace) re-merges them into a single run (no visible text\n#    change, only a run-structure normalization).\n$colon = $d.Content\n$find = $colon.Find\n$find.ClearFormatting()\n$find.Text = \"but there\u2019s a problem:\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"but there\u2019s a problem:\"\n$wdReplaceOne = 1\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, $wdReplaceOne) | Out-Null\n\n# 2) Insert a new clause right after \"...You bring no silver at all\" and\n#    before the following \". How would you reach the temple...\" sentence.\n$anchor = $d.Content\n$anchor.Find.Execute(\"You bring no silver at all\") | Out-Null\n$anchor.Collapse(0)  # wdCollapseEnd\n$newText = \", and you can\u2019t walk the path you\u2019ve trodden, though you can cross path\"\n$anchor.InsertAfter($newText)\n\n# Force the newly inserted clause into its own run (distinct from the\n# surrounding text it was merged with) by toggling a character property\n# on it and back off again. $anchor now spans exactly the inserted text.\n$anchor.Font.Bold = 1\n$anchor.Font.Bold = 0\n"}
